# Updates the cryptos price/volume table to the latest scraped values
# (commit: "Updated cryptos list on Thu Jun 22 07:21:03 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage ("@") before writing, since these columns hold
# human-formatted strings (e.g. "0.7160", "30.158.71") that must not be
# auto-coerced to numbers by the COM Value setter.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.158.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5110"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.05"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2963"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06811"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.909.31"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.26"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07364"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6915"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008365"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +13.25%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.150.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.157.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.813"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.710"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.135"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.79"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.998"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.227"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08821"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.010"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05062"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.148"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7160"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.810"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.279"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01689"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.170"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4293"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.29%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.618"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1277"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05731"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.396"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3804"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.87%  "
